$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Highlight rows 8-13 (the "O" / output rows) with a yellow fill and make
#    sure the font used on those rows is re-asserted (Calibri 11) - this
#    mirrors the author selecting A8:S13 and applying a yellow background.
#    The block is processed as five contiguous rectangles that each share a
#    uniform pre-existing style so the engine does not fragment the style
#    table unnecessarily.
# ---------------------------------------------------------------------------
$highlightBlocks = @("A8:C13", "L8:O13", "R8:S13", "D8:K13", "P8:Q13")
foreach ($block in $highlightBlocks) {
    $rng = $ws.Range($block)
    $rng.Font.Name = "Calibri"
    $rng.Interior.Color = 65535
}

# ---------------------------------------------------------------------------
# 2. Add the new T3/U3 cells: T3 just gets the same yellow fill, U3 gets a
#    (broken / not-yet-implemented) formula that evaluates to a #NAME? error.
# ---------------------------------------------------------------------------
$ws.Range("T3").Interior.Color = 65535
$ws.Range("U3").Formula = "= not yet implemented in demonstrator"

# ---------------------------------------------------------------------------
# 3. Update the view: zoom to 85%, scroll so column J is near the left edge,
#    and move the selection to the new U3 cell.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 85
$win.ScrollColumn = 10
$win.ScrollRow = 1

$ws.Range("U3").Select()
